# Split the single run-on "Bibliografia" paragraph into separate
# reference lines, joined by manual line breaks (<w:br/>), matching
# the target diff. We rebuild the paragraph's run content via
# InsertXML so we get byte-exact control over which <w:t> elements
# carry xml:space="preserve" (only the two entries that genuinely end
# with a trailing space need it).

$d = $word.ActiveDocument

# Locate the bibliography paragraph (the one whose text starts with
# "BANZATO").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("BANZATO")) {
        $target = $cand
        break
    }
}

$insertPoint = $target.Range.Start

$runInner = '<w:t>BANZATO, Eduardo et al. Atualidades na armazenagem. São Paulo: IMAM, 2003.</w:t><w:br/><w:t>BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. São Paulo, Edgar Blücher, 1977.</w:t><w:br/><w:t xml:space="preserve">GURGEL, F.A.C. Administração de recursos materiais e patrimoniais. 2a. Edição. São Paulo. Editora Cengage. 2013. </w:t><w:br/><w:t>FRANCISCHINI, P.G.; VALLE, C.E. Implantação de Indústrias. Rio de Janeiro, LTC Editora, 1975.</w:t><w:br/><w:t>LEE, Q et al. Projeto de Instalações e Locais de Trabalho. São Paulo: IMAM, 1998.</w:t><w:br/><w:t>MOURA, Reinaldo Aparecido. Sistemas e técnicas de movimentação e armazenagem de materiais. IMAM, 2012.</w:t><w:br/><w:t>NEWMANN, C.; SCALICE, R.K. Projeto de Fábrica e Layout. Rio de Janeiro, Elsevier, 2015.</w:t><w:br/><w:t xml:space="preserve">Müther, R. Planejamento do Layout: Sistema SLP. São Paulo, Edgard Blücher, 1978. </w:t><w:br/><w:t>SLACK, Nigel et al. Administração da produção. São Paulo: Atlas, 8ª ed. 2018.</w:t><w:br/><w:t>TOMPKINS, James A. et al. Planejamento de instalações. Editora LTC:, 2013.</w:t>'

$xmlSnippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Insert the replacement content right before the existing (old)
# text. This merges into the current paragraph rather than creating a
# new one, since the insertion point is collapsed and sits at the
# very start of the paragraph's existing content.
$rngInsert = $d.Range($insertPoint, $insertPoint)
$rngInsert.InsertXML($xmlSnippet)

# Now remove the old run-on text that follows what we just inserted.
# (Use LastIndexOf because our freshly-inserted text also starts with
# "BANZATO" -- the old copy is the later occurrence.)
$p2 = $target
$fullText = $p2.Range.Text
$oldStartIdx = $fullText.LastIndexOf("BANZATO")

$delStart = $p2.Range.Start + $oldStartIdx
$delEnd = $p2.Range.End - 1    # stop before the paragraph mark
$rngDel = $d.Range($delStart, $delEnd)
$rngDel.Delete()
